$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36/37 swap: Celestia <-> FirstDigitalUSD, with new ranking & values
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.96"
$ws.Range("E37").Value = "  -1.26%  "

# Price / volume updates
$ws.Range("D2").Value = "51.490.18"
$ws.Range("E2").Value = "  +6.04%  "
$ws.Range("D3").Value = "2.739.43"
$ws.Range("E3").Value = "  +4.63%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "116.57"
$ws.Range("E5").Value = "  +6.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "331.03"
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.531"
$ws.Range("E7").Value = "  +2.71%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.573"
$ws.Range("E9").Value = "  +6.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.28"
$ws.Range("E10").Value = "  +5.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.92"
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0826"
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.128"
$ws.Range("E13").Value = "  +2.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.57"
$ws.Range("E14").Value = "  +5.62%  "
$ws.Range("D15").Value = "3.182.07"
$ws.Range("E15").Value = "  +5.29%  "
$ws.Range("D16").Value = "2.765.58"
$ws.Range("E16").Value = "  +5.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.879"
$ws.Range("E17").Value = "  +2.93%  "
$ws.Range("D18").Value = "51.550.48"
$ws.Range("E18").Value = "  +6.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.50"
$ws.Range("E19").Value = "  +5.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.99"
$ws.Range("E20").Value = "  +3.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.80"
$ws.Range("E21").Value = "  +2.66%  "
$ws.Range("D22").Value = "0.0₃0958"
$ws.Range("E22").Value = "  +2.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "277.69"
$ws.Range("E23").Value = "  +3.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.46"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  +4.44%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.66"
$ws.Range("E26").Value = "  +2.53%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.13"
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.18"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.139"
$ws.Range("E31").Value = "  +2.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.69"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "50.23"
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.55"
$ws.Range("E34").Value = "  +3.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0817"
$ws.Range("E35").Value = "  +3.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.08"
$ws.Range("E38").Value = "  +3.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.92"
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.18"
$ws.Range("E40").Value = "  +2.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "130.82"
$ws.Range("E41").Value = "  +4.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "23.16"
$ws.Range("E42").Value = "  +4.11%  "
$ws.Range("E43").Value = "  +11.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.113"
$ws.Range("E44").Value = "  +2.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.26"
$ws.Range("E45").Value = "  +5.10%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.37"
$ws.Range("E46").Value = "  +14.15%  "
$ws.Range("D47").Value = "2.103.26"
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.31"
$ws.Range("E48").Value = "  +4.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.24"
$ws.Range("E49").Value = "  +2.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.53"
$ws.Range("E50").Value = "  +7.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.93"
$ws.Range("E51").Value = "  +0.52%  "
